# Fruta / hortaliza, semanal
# Insert a new weekly record at row 50 (pushing existing rows 50-149 down to 51-150)
# and populate the new row with the latest price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 50; this shifts rows 50:149 down to 51:150
$ws.Rows.Item(50).Insert()

# Populate the new row 50 with the latest weekly data
$ws.Range("A50").Value = 8
$ws.Range("B50").Value = "Terminal La Palmera de La Serena"
$ws.Range("C50").Value = "Coquimbo"
$ws.Range("D50").Value = 44967
$ws.Range("E50").Value = 4
$ws.Range("F50").Value = 100112028
$ws.Range("G50").Value = "Sandia"
$ws.Range("H50").Value = "Sin especificar"
$ws.Range("I50").Value = "Extra"
$ws.Range("J50").Value = 2000
$ws.Range("K50").Value = 3300
$ws.Range("L50").Value = 3500
$ws.Range("M50").Value = 3400
$ws.Range("N50").Value = "`$/unidad"
$ws.Range("O50").Value = "Región de O'Higgins"
$ws.Range("P50").Value = 3400
$ws.Range("Q50").Value = 1
$ws.Range("R50").Value = "Hortaliza"
